$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column L
$ws.Range("L1").Value = "STD/MEAN"

# Fill in STD/MEAN formula for every data row (2 through 61), referencing the
# "youngs_moduli_korrigiert" table via structured references.
for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $cell.Formula = "=youngs_moduli_korrigiert[[#This Row],[E_Std]]/youngs_moduli_korrigiert[[#This Row],[E_Mean_corrected]]"
}

# Match the saved selection state from the diff
$ws.Range("L2:L38").Select()
